# Apply crypto price/volume updates to the active worksheet
# (mirrors the GitHub Actions "Updated cryptos list" commit).
#
# Column D ("Price") holds free-form text like "59.194.42" or
# "2.991.00" that LOOKS numeric to Excel's input parser, so a plain
# Value assignment would silently convert it to a real number. We
# force the General-formatted cell to Text just for the assignment,
# then flip it back to the built-in "Normal" style so no stray
# number-format / style residue is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '59.194.42'
$ws.Range('E2').Value = '  +3.14%  '
Set-TextValue 'D3' '2.991.00'
$ws.Range('E3').Value = '  +2.10%  '
$ws.Range('E4').Value = '  -0.18%  '
Set-TextValue 'D5' '562.45'
$ws.Range('E5').Value = '  +2.10%  '
Set-TextValue 'D6' '137.81'
$ws.Range('E6').Value = '  +6.27%  '
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('E8').Value = '  +1.68%  '
Set-TextValue 'D9' '2.983.17'
$ws.Range('E9').Value = '  +2.19%  '
Set-TextValue 'D10' '0.133'
$ws.Range('E10').Value = '  +3.92%  '
$ws.Range('E11').Value = '  +6.59%  '
$ws.Range('E12').Value = '  +2.96%  '
Set-TextValue 'D13' '0.0000230'
$ws.Range('E13').Value = '  +4.32%  '
Set-TextValue 'D14' '33.64'
$ws.Range('E14').Value = '  +3.90%  '
$ws.Range('E15').Value = '  +2.44%  '
Set-TextValue 'D16' '3.485.19'
$ws.Range('E16').Value = '  +2.06%  '
Set-TextValue 'D17' '7.22'
$ws.Range('E17').Value = '  +7.58%  '
Set-TextValue 'D18' '2.988.99'
$ws.Range('E18').Value = '  +2.05%  '
Set-TextValue 'D19' '59.144.09'
$ws.Range('E19').Value = '  +2.67%  '
Set-TextValue 'D20' '429.15'
$ws.Range('E20').Value = '  +3.51%  '
$ws.Range('E21').Value = '  +5.01%  '
$ws.Range('E22').Value = '  +6.15%  '
$ws.Range('E23').Value = '  +2.78%  '
Set-TextValue 'D24' '13.30'
$ws.Range('E24').Value = '  +3.38%  '
Set-TextValue 'D25' '80.77'
$ws.Range('E25').Value = '  +2.60%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('E28').Value = '  +9.23%  '
$ws.Range('E29').Value = '  +2.97%  '
Set-TextValue 'D30' '7.77'
$ws.Range('E30').Value = '  +4.44%  '
Set-TextValue 'D31' '25.71'
$ws.Range('E31').Value = '  +3.36%  '
$ws.Range('E32').Value = '  -0.39%  '
Set-TextValue 'D33' '0.0988'
$ws.Range('E33').Value = '  -4.24%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D34' '5.92'
$ws.Range('E34').Value = '  +5.87%  '
$ws.Range('B35').Value = 'Mantle'
$ws.Range('C35').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D35' '0.991'
$ws.Range('E35').Value = '  +6.70%  '
Set-TextValue 'D36' '0.0₃0765'
$ws.Range('E36').Value = '  +13.50%  '
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('E38').Value = '  +1.06%  '
$ws.Range('E39').Value = '  +3.86%  '
Set-TextValue 'D40' '2.71'
$ws.Range('E40').Value = '  +6.54%  '
Set-TextValue 'D41' '400.20'
$ws.Range('E41').Value = '  +6.39%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D42' '2.756.06'
$ws.Range('E42').Value = '  +4.76%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D43' '0.0351'
$ws.Range('E43').Value = '  +1.10%  '
$ws.Range('E44').Value = '  +0.53%  '
Set-TextValue 'D45' '0.252'
$ws.Range('E45').Value = '  +5.94%  '
$ws.Range('E46').Value = '  +0.01%  '
Set-TextValue 'D47' '34.88'
$ws.Range('E47').Value = '  +22.62%  '
Set-TextValue 'D48' '121.32'
$ws.Range('E48').Value = '  +0.41%  '
$ws.Range('E49').Value = '  +1.86%  '
$ws.Range('E50').Value = '  +1.48%  '
Set-TextValue 'D51' '23.40'
$ws.Range('E51').Value = '  +1.13%  '
